$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style from an existing header cell (AC1) to the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

for ($r = 2; $r -le 63; $r++) {
    $ws.Cells.Item($r, 30).Value = 78   # AD
    $ws.Cells.Item($r, 31).Value = 84   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
